$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph "Elaborazione - Iterazione 1": give the paragraph mark
#    and the run an explicit font size of 16pt (sz/szCs = 32 half-points).
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.Font.Size = 16
$title.Range.Font.SizeBi = 16

# ---------------------------------------------------------------------------
# 2) Table cell "e stato restituito un Messaggio di Verifica dal Sistema.":
#    capitalize the leading "e" -> "E" (with accent) and split it into its
#    own run, leaving the remainder of the sentence as a second run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("è stato restituito un Messaggio di Verifica dal Sistema.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start
    $end = $rng.End

    # Force a run boundary right after the first character by toggling a
    # character-formatting property on just that character, then revert it
    # so the visible formatting is unchanged but the run has been split.
    $firstChar = $d.Range($start, $start + 1)
    $firstChar.Bold = 1
    $firstChar.Text = "È"
    $firstChar.Bold = 0
}
